# Update phone number formatting for three providers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D9").Value = "(844) 800-6020"
$ws.Range("D15").Value = "(888) 508-3028-101"
$ws.Range("D22").Value = "(506) 536-4565"
